$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): extend sequence with P1=14, Q1=15 ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
# Match header style (bold / bordered) used by the other header cells in row 1
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows 2-25: swap I<->K and M<->O values, then add P=2 / Q=2 ---
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # column I
    $kVal = $ws.Cells.Item($r, 11).Value2  # column K
    $ws.Cells.Item($r, 9).Value = $kVal
    $ws.Cells.Item($r, 11).Value = $iVal

    $mVal = $ws.Cells.Item($r, 13).Value2  # column M
    $oVal = $ws.Cells.Item($r, 15).Value2  # column O
    $ws.Cells.Item($r, 13).Value = $oVal
    $ws.Cells.Item($r, 15).Value = $mVal

    $ws.Cells.Item($r, 16).Value = 2      # column P
    $ws.Cells.Item($r, 17).Value = 2      # column Q
}
